$p = $ppt.ActivePresentation

# --- 1) Update the "Date Placeholder" field text (1/29/2025 -> 2/17/2025)
#        on the slide master and on every slide layout. ---
function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "1/29/2025") {
                $tr.Text = "2/17/2025"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# --- 2) Slide 4: split the last sentence of the "speed" paragraph into
#        its own (underlined) run: "As long as you respond within 2.5s, "
#        + underlined "speed is NOT critical here" ---
$slide4 = $p.Slides.Item(4)
$shape = $slide4.Shapes.Item(1)
$textRange = $shape.TextFrame.TextRange

$target = "As long as you respond within 2.5s, speed is NOT critical here"
$paraCount = $textRange.Paragraphs().Count
for ($pi = 1; $pi -le $paraCount; $pi++) {
    $para = $textRange.Paragraphs($pi, 1)
    if ($para.Text.TrimEnd() -eq $target) {
        $prefix = "As long as you respond within 2.5s, "
        $splitLen = $prefix.Length
        $tail = $textRange.Characters($para.Start + $splitLen, $para.Length - $splitLen)
        $tail.Font.Underline = $true
    }
}
